$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: player/killer/time cleared, points reset to 0
$ws.Range("B3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = 0
$ws.Range("F3").ClearContents()

# Row 4
$ws.Range("B4").Value = "Didier"
$ws.Range("D4").Value = "Eric"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "10:05"

# Row 5
$ws.Range("B5").Value = "Baptiste"
$ws.Range("D5").Value = "Didier"
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = "10:05"

# Row 6 (player unchanged, still Côme)
$ws.Range("D6").Value = "Didier"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = "10:05"

# Row 7
$ws.Range("B7").Value = "Sylvie P"
$ws.Range("D7").Value = "Didier"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "10:04"

# Row 8
$ws.Range("B8").Value = "Béa"
$ws.Range("D8").Value = "Didier"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "10:04"

# Row 9
$ws.Range("B9").Value = "Anne-Lise"
$ws.Range("D9").Value = "Didier"
$ws.Range("F9").Value = "10:04"

# Row 10
$ws.Range("B10").Value = "Mathieu"
$ws.Range("D10").Value = "Eric"
$ws.Range("F10").Value = "10:04"

# Row 11
$ws.Range("B11").Value = "Jean Rob"
$ws.Range("D11").Value = "Eric"
$ws.Range("F11").Value = "10:04"

# Row 12
$ws.Range("B12").Value = "Hugo D"
$ws.Range("D12").Value = "Jean Rob"
$ws.Range("F12").Value = "10:03"
